$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C5: "En proceso" -> "Leido"
$ws.Range("C5").Value = "Leido"

# C6: new value "Leido"
$ws.Range("C6").Value = "Leido"

# C8: new value "enProceso"
$ws.Range("C8").Value = "enProceso"

# Update selection to C8
$ws.Range("C8").Select()
